$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title: "PVV-politicus Martin Bosma verrast door taartaanval"
#    -> "PVV’er Martin Bosma verrast door taartaanval"
# ------------------------------------------------------------------
$d.Content.Find.Execute("PVV-politicus Martin Bosma verrast door taartaanval", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "PVV’er Martin Bosma verrast door taartaanval", 2) | Out-Null

# ------------------------------------------------------------------
# 2. "... had Martin Bosma wellicht niet verwacht ..."
#    -> "... had Martin Bosma waarschijnlijk niet verwacht ..."
# ------------------------------------------------------------------
$d.Content.Find.Execute(" wellicht niet verwacht", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " waarschijnlijk niet verwacht", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Rewrite the "Taartgooier weggekomen" body paragraph completely.
#    (Old paragraph 8, 1-indexed.) Replace the paragraph's text
#    (excluding its trailing paragraph mark) with the new wording,
#    dropping the bold "Het incident..." run and the trailing
#    " De PVV’er" fragment.
# ------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Het is niet bekend of de jongeren")) {
        $targetPara = $p
        break
    }
}
if ($targetPara -ne $null) {
    $r = $targetPara.Range
    $r.End = $r.End - 1
    $r.Text = "De dader is weggekomen en het is niet bekend of hij een student aan Leiden is, verteld de woordvoerder Faculteit Rechten Universiteit Leiden. Wel loopt er een onderzoek Het incident diende als een opmerkelijke uiting van protest tegen de recente PVV-overwinning. "
}

# ------------------------------------------------------------------
# 4. Merge the "Staan open voor opnieuw ..." paragraph with the
#    "Universiteit wil alle partijen een stem geven, " and
#    "Bedoeling dat er nieuwe datum wordt geprikt, door
#    verkiezingen onduidelijk" paragraphs, dropping the entire
#    "Niet perse inzetten op strengere beveiliging ..." paragraph
#    and the blank paragraph that followed it.
# ------------------------------------------------------------------
# Find paragraph indices by scanning (robust against any prior shifts).
function Get-ParaIndexStartingWith($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

$iStaan = Get-ParaIndexStartingWith("Staan open voor opnieuw")
$iNietPerse = Get-ParaIndexStartingWith("Niet perse inzetten")
$iBedoeling = Get-ParaIndexStartingWith("Bedoeling dat er nieuwe datum")

if ($iStaan -gt 0 -and $iBedoeling -gt 0) {
    $pStaan = $d.Paragraphs.Item($iStaan)
    $rStaan = $pStaan.Range
    $rStaan.End = $rStaan.End - 1
    $rStaan.Text = "Staan open voor opnieuw martin uit te nodigen Hoe gaat Leiden dit soort dingen voorkomen: Ze gaan beveiliging optrekken, LU kaarten actief Universiteit wil alle partijen een stem geven, Bedoeling dat er nieuwe datum wordt geprikt, door verkiezingen onduidelijk"

    # Delete everything from the (now obsolete) "Niet perse ..." paragraph
    # through the end of the (now empty) "Bedoeling ..." paragraph, i.e.
    # through the blank paragraph right after it.
    $iNietPerse2 = Get-ParaIndexStartingWith("Niet perse inzetten")
    if ($iNietPerse2 -gt 0) {
        $startDel = $d.Paragraphs.Item($iNietPerse2).Range.Start
        # walk forward from $iNietPerse2 until we reach the paragraph right
        # before "Staart:" — that is, consume the "Niet perse..." paragraph,
        # the now-empty former "Universiteit wil..." paragraph, the
        # now-empty former "Bedoeling..." paragraph, and the blank paragraph.
        $iStaart = Get-ParaIndexStartingWith("Staart:")
        $endDel = $d.Paragraphs.Item($iStaart - 1).Range.End
        $delRng = $d.Range($startDel, $endDel)
        $delRng.Delete() | Out-Null
    }
}
